$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking price strings are not
# auto-converted to numbers by Excel, matching the original inline-string cells.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "27.853.31"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.768.19"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "327.54"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "0.4483"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("D8").Value = "0.3550"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").Value = "0.07452"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "42.12"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").Value = "1.095"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "0.9996"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "20.84"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "6.026"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "7.200"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "1.768.89"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "93.11"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "0.00001057"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "0.06434"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("D22").Value = "5.773"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "27.901.36"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "162.21"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").Value = "20.22"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "1.969.12"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "2.160"
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("D30").Value = "125.23"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "1.094"
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").Value = "0.09153"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "3.653"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "5.575"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "11.89"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "0.02291"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "0.06100"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("D38").Value = "0.2095"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "0.6291"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "7.925"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").Value = "13.23"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "3.740"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "0.5857"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "122.39"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "1.946"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "0.06909"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "72.74"
$ws.Range("E51").Value = "  +0.98%  "

# Restore default (un-styled) cell formatting so no stray style index
# is left attached to these cells, matching the original workbook.
$rng.Style = "Normal"
